# Weekly data update: insert the newest "Poroto granado" price record
# for "Terminal La Palmera de La Serena" at the top of the data block
# (row 41), pushing the previous entries down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 41, shifting rows 41:45 down to 42:46
$ws.Rows.Item(41).Insert()

# Populate the new row 41 with the latest weekly record
$ws.Cells.Item(41,1).Value = 8
$ws.Cells.Item(41,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(41,3).Value = "Coquimbo"
$ws.Cells.Item(41,4).Value = 44505
$ws.Cells.Item(41,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(41,5).Value = 4
$ws.Cells.Item(41,6).Value = 100112030
$ws.Cells.Item(41,7).Value = "Poroto granado"
$ws.Cells.Item(41,8).Value = "Sin especificar"
$ws.Cells.Item(41,9).Value = "Primera"
$ws.Cells.Item(41,10).Value = 300
$ws.Cells.Item(41,11).Value = 37000
$ws.Cells.Item(41,12).Value = 38000
$ws.Cells.Item(41,13).Value = 37500
$ws.Cells.Item(41,14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(41,15).Value = "Perú"
$ws.Cells.Item(41,16).Value = 1500
$ws.Cells.Item(41,17).Value = 25
$ws.Cells.Item(41,18).Value = "Hortaliza"
